$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: order quantity + total amount change
$ws.Range("B2").Value = 5
$ws.Range("G2").Value = 132389000

# Row 3: now refers to Nguyễn Thị Diệu Mỵ's order (was a second Đỗ Minh Tâm row)
$ws.Range("B3").Value = 8
$ws.Range("C3").Value = "Nguyễn Thị Diệu Mỵ"
$ws.Range("E3").Value = "10:12:05 11/06/2024"
$ws.Range("G3").Value = 8994000

# Row 4: now refers to vanh123's first order
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = "vanh123"
$ws.Range("E4").Value = "01:41:26 12/06/2024"
$ws.Range("G4").Value = 83098000

# Row 5: now refers to vanh123's second order
$ws.Range("B5").Value = 10
$ws.Range("C5").Value = "vanh123"
$ws.Range("E5").Value = "01:45:27 12/06/2024"
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 67799000

# Remove the old rows 6 and 7 (their data now folded into rows above / removed)
$ws.Rows("6:7").Delete()
